# edit.ps1
# Applies the Diario_2021_04_15.docx edit:
#   1. Merge the "Sistemata la possibilita` ... cartella." run with the
#      trailing space-only run into a single run (and drop the _GoBack
#      bookmark that used to sit between them).
#   2. Insert a new empty "List Paragraph" (indent 1068 twips) carrying
#      the (relocated) _GoBack bookmark right before the "Per Thaisa:"
#      paragraph, and move the <w:lastRenderedPageBreak/> marker from
#      the "Recuperare le 4 ore della mattinata" run onto the
#      "Per Thaisa:" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: "Sistemata la possibilita` della destinazione singola..."
# Find/Execute spans the run boundary between the sentence and the
# trailing " " run (which also removes the _GoBack bookmark that used
# to live between those two runs), merging everything into one run.
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Sistemata la possibilit") {
        $found = $true
        $target = $para
    }
}
if ($found) {
    $rng = $target.Range
    $rng.Find.ClearFormatting()
    [void]$rng.Find.Execute("cartella. ", $false, $false, $false, $false, $false, $true, 1, $false, "cartella. ", 2)
}

# ---------------------------------------------------------------------
# Part 2: locate the "Sistemare errore permessi..." paragraph (the one
# right before "Per Thaisa:") and splice in a brand-new empty paragraph
# after it, carrying the relocated _GoBack bookmark.
# ---------------------------------------------------------------------
$sepIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Sistemare errore permessi") {
        $sepIdx = $i
    }
}

if ($sepIdx -ge 0) {
    $sepPara = $d.Paragraphs.Item($sepIdx)
    $insertAt = $d.Range($sepPara.Range.End, $sepPara.Range.End)
    $bookmarkParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:ind w:left="1068"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$insertAt.InsertXML($bookmarkParaXml)
}

# ---------------------------------------------------------------------
# Part 3: move <w:lastRenderedPageBreak/> from the "Recuperare le 4 ore
# della mattinata" run onto the "Per Thaisa:" run, rewriting each
# paragraph surgically (keeping their original rsid attributes) so only
# that element's position changes.
# ---------------------------------------------------------------------
$thaisaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^\s*Per Thaisa:\s*$") {
        $thaisaIdx = $i
    }
}
if ($thaisaIdx -ge 0) {
    $thaisaPara = $d.Paragraphs.Item($thaisaIdx)
    $thaisaRng = $d.Range($thaisaPara.Range.Start, $thaisaPara.Range.End)
    $thaisaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005A1496" w:rsidRDefault="005A1496" w:rsidP="005A1496"><w:r><w:lastRenderedPageBreak/><w:t>Per Thaisa:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$thaisaRng.InsertXML($thaisaXml)
}

$recIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Recuperare le 4 ore della mattinata") {
        $recIdx = $i
    }
}
if ($recIdx -ge 0) {
    $recPara = $d.Paragraphs.Item($recIdx)
    $recRng = $d.Range($recPara.Range.Start, $recPara.Range.End)
    $recXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005A1496" w:rsidRPr="005A1496" w:rsidRDefault="005A1496" w:rsidP="002A7EB0"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Recuperare le 4 ore della mattinata</w:t></w:r><w:r w:rsidR="004536ED"><w:t xml:space="preserve"> del 25/03</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$recRng.InsertXML($recXml)
}

Write-Host "Edit applied."
